$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 17.70643966666667
$ws.Range("H2").Value = 53.119319
$ws.Range("I2").Value = 0.4380235920947999
$ws.Range("J2").Value = 0.4380235920947999
$ws.Range("O2").Value = 0.7091726973716084
$ws.Range("P2").Value = 0.7091726973716084
$ws.Range("Q2").Value = 30.23882157687111
$ws.Range("R2").Value = 272.14939419184
$ws.Range("S2").Value = 0.3106343723182703
$ws.Range("T2").Value = 0.3106343723182703

# Row 3
$ws.Range("G3").Value = 17.70643966666667
$ws.Range("H3").Value = 53.119319
$ws.Range("I3").Value = 0.4380235920947999
$ws.Range("J3").Value = 0.4380235920947999
$ws.Range("M3").Value = 0.7003526666666667
$ws.Range("N3").Value = 2.101058
$ws.Range("O3").Value = 0.2908273026283917
$ws.Range("P3").Value = 0.2908273026283917
$ws.Range("Q3").Value = 12.40075223772245
$ws.Range("R3").Value = 111.606770139502
$ws.Range("S3").Value = 0.1273892197765296
$ws.Range("T3").Value = 0.1273892197765296

# Row 4
$ws.Range("G4").Value = 1.617245333333334
$ws.Range("H4").Value = 4.851736000000001
$ws.Range("I4").Value = 0.04000756919748267
$ws.Range("J4").Value = 0.04000756919748267
$ws.Range("O4").Value = 0.7091726973716084
$ws.Range("P4").Value = 0.7091726973716084
$ws.Range("Q4").Value = 2.761910016995556
$ws.Range("S4").Value = 0.02837227576306006
$ws.Range("T4").Value = 0.02837227576306006

# Row 5
$ws.Range("G5").Value = 1.617245333333334
$ws.Range("H5").Value = 4.851736000000001
$ws.Range("I5").Value = 0.04000756919748267
$ws.Range("J5").Value = 0.04000756919748267
$ws.Range("M5").Value = 0.7003526666666667
$ws.Range("N5").Value = 2.101058
$ws.Range("O5").Value = 0.2908273026283917
$ws.Range("P5").Value = 0.2908273026283917
$ws.Range("Q5").Value = 1.132642081854223
$ws.Range("R5").Value = 10.193778736688
$ws.Range("S5").Value = 0.01163529343442261
$ws.Range("T5").Value = 0.01163529343442261

# Row 6
$ws.Range("G6").Value = 21.099799
$ws.Range("H6").Value = 63.299397
$ws.Range("I6").Value = 0.5219688387077175
$ws.Range("J6").Value = 0.5219688387077175
$ws.Range("O6").Value = 0.7091726973716084
$ws.Range("P6").Value = 0.7091726973716084
$ws.Range("Q6").Value = 36.03395540154666
$ws.Range("R6").Value = 324.30559861392
$ws.Range("S6").Value = 0.370166049290278
$ws.Range("T6").Value = 0.370166049290278

# Row 7
$ws.Range("G7").Value = 21.099799
$ws.Range("H7").Value = 63.299397
$ws.Range("I7").Value = 0.5219688387077175
$ws.Range("J7").Value = 0.5219688387077175
$ws.Range("M7").Value = 0.7003526666666667
$ws.Range("N7").Value = 2.101058
$ws.Range("O7").Value = 0.2908273026283917
$ws.Range("P7").Value = 0.2908273026283917
$ws.Range("Q7").Value = 14.77730049578067
$ws.Range("R7").Value = 132.995704462026
$ws.Range("S7").Value = 0.1518027894174395
$ws.Range("T7").Value = 0.1518027894174395
